# Apply cryptos list update (commit: "Updated cryptos list on Sat May 27 07:05:06 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 20 & 21 swapped places (Avalanche <-> WrappedBTC) plus value refresh ---
$ws.Cells.Item(20, 2).Value = "WrappedBTC"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(20, 4).Value = "26.948.78"
$ws.Cells.Item(20, 5).Value = "  +1.18%  "

$ws.Cells.Item(21, 2).Value = "Avalanche"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(21, 4).Value = "'14.56"
$ws.Cells.Item(21, 5).Value = "  +2.89%  "

# --- Per-row Price (D) / Volume(1h) (E) refresh ---
$ws.Cells.Item(2, 4).Value = "26.922.61"
$ws.Cells.Item(2, 5).Value = "  +1.22%  "
$ws.Cells.Item(3, 4).Value = "1.845.39"
$ws.Cells.Item(3, 5).Value = "  +1.23%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).Value = "'309.42"
$ws.Cells.Item(5, 5).Value = "  +0.33%  "
$ws.Cells.Item(6, 5).Value = "  +0.00%  "
$ws.Cells.Item(7, 4).Value = "'0.4758"
$ws.Cells.Item(7, 5).Value = "  +2.62%  "
$ws.Cells.Item(8, 4).Value = "'0.3668"
$ws.Cells.Item(8, 5).Value = "  +1.81%  "
$ws.Cells.Item(9, 4).Value = "'0.07201"
$ws.Cells.Item(9, 5).Value = "  +0.95%  "
$ws.Cells.Item(10, 4).Value = "'0.9268"
$ws.Cells.Item(10, 5).Value = "  +2.92%  "
$ws.Cells.Item(11, 4).Value = "'19.74"
$ws.Cells.Item(11, 5).Value = "  +1.92%  "
$ws.Cells.Item(12, 4).Value = "'0.07692"
$ws.Cells.Item(12, 5).Value = "  -1.01%  "
$ws.Cells.Item(13, 4).Value = "1.829.78"
$ws.Cells.Item(13, 5).Value = "  +0.64%  "
$ws.Cells.Item(14, 4).Value = "'5.318"
$ws.Cells.Item(14, 5).Value = "  +1.02%  "
$ws.Cells.Item(15, 4).Value = "'6.405"
$ws.Cells.Item(15, 5).Value = "  +1.39%  "
$ws.Cells.Item(16, 4).Value = "'88.75"
$ws.Cells.Item(16, 5).Value = "  +1.52%  "
$ws.Cells.Item(17, 4).Value = "'1.010"
$ws.Cells.Item(17, 5).Value = "  +0.11%  "
$ws.Cells.Item(18, 4).Value = "'0.000008641"
$ws.Cells.Item(18, 5).Value = "  +0.96%  "
$ws.Cells.Item(19, 5).Value = "  +0.04%  "
$ws.Cells.Item(22, 4).Value = "'5.049"
$ws.Cells.Item(22, 5).Value = "  +0.69%  "
$ws.Cells.Item(23, 4).Value = "'10.63"
$ws.Cells.Item(23, 5).Value = "  +0.86%  "
$ws.Cells.Item(24, 4).Value = "'1.922"
$ws.Cells.Item(24, 5).Value = "  +0.26%  "
$ws.Cells.Item(25, 4).Value = "'152.40"
$ws.Cells.Item(25, 5).Value = "  +0.31%  "
$ws.Cells.Item(26, 4).Value = "'18.15"
$ws.Cells.Item(26, 5).Value = "  +1.36%  "
$ws.Cells.Item(27, 4).Value = "'1.999"
$ws.Cells.Item(27, 5).Value = "  +1.06%  "
$ws.Cells.Item(28, 4).Value = "'114.24"
$ws.Cells.Item(28, 5).Value = "  +0.40%  "
$ws.Cells.Item(29, 4).Value = "'4.930"
$ws.Cells.Item(29, 5).Value = "  +2.62%  "
$ws.Cells.Item(30, 4).Value = "'0.08879"
$ws.Cells.Item(30, 5).Value = "  +0.83%  "
$ws.Cells.Item(31, 4).Value = "'3.307"
$ws.Cells.Item(31, 5).Value = "  +5.46%  "
$ws.Cells.Item(32, 4).Value = "'1.175"
$ws.Cells.Item(32, 5).Value = "  +3.37%  "
$ws.Cells.Item(33, 4).Value = "'0.7485"
$ws.Cells.Item(33, 5).Value = "  +2.23%  "
$ws.Cells.Item(34, 4).Value = "'4.483"
$ws.Cells.Item(34, 5).Value = "  +1.05%  "
$ws.Cells.Item(35, 4).Value = "'2.743"
$ws.Cells.Item(35, 5).Value = "  +0.55%  "
$ws.Cells.Item(36, 4).Value = "'1.096"
$ws.Cells.Item(36, 5).Value = "  +2.09%  "
$ws.Cells.Item(37, 4).Value = "'0.01957"
$ws.Cells.Item(37, 5).Value = "  +1.84%  "
$ws.Cells.Item(38, 4).Value = "'0.05269"
$ws.Cells.Item(38, 5).Value = "  +3.04%  "
$ws.Cells.Item(39, 5).Value = "  +1.74%  "
$ws.Cells.Item(40, 5).Value = "  +2.92%  "
$ws.Cells.Item(41, 4).Value = "'6.982"
$ws.Cells.Item(41, 5).Value = "  +1.34%  "
$ws.Cells.Item(42, 4).Value = "'0.1511"
$ws.Cells.Item(42, 5).Value = "  +1.34%  "
$ws.Cells.Item(43, 4).Value = "'8.227"
$ws.Cells.Item(43, 5).Value = "  +2.91%  "
$ws.Cells.Item(44, 4).Value = "'10.55"
$ws.Cells.Item(44, 5).Value = "  +5.87%  "
$ws.Cells.Item(45, 5).Value = "  +1.61%  "
$ws.Cells.Item(46, 5).Value = "  +0.01%  "
$ws.Cells.Item(47, 4).Value = "'101.55"
$ws.Cells.Item(47, 5).Value = "  +3.40%  "
$ws.Cells.Item(48, 4).Value = "'1.605"
$ws.Cells.Item(48, 5).Value = "  +3.06%  "
$ws.Cells.Item(49, 4).Value = "'66.10"
$ws.Cells.Item(49, 5).Value = "  +3.82%  "
$ws.Cells.Item(50, 4).Value = "'0.06023"
$ws.Cells.Item(50, 5).Value = "  +0.72%  "
$ws.Cells.Item(51, 4).Value = "'0.8852"
$ws.Cells.Item(51, 5).Value = "  +4.01%  "
